# Updated cryptos list on Wed Dec 27 21:55:21 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $value) {
    # Force the cell to keep its value as text (the "Price" column holds
    # numeric-looking strings like "43.18" or "2.714.84" that must stay
    # literal text, not be re-interpreted/rounded as numbers), then drop
    # back to the Normal style so no extra formatting is introduced.
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "43.425.42"
$ws.Range("E2").Value = "  +2.60%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.357.71"
$ws.Range("E3").Value = "  +6.08%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.16%  "

# Row 5 - BNB
Set-TextValue "D5" "311.38"
$ws.Range("E5").Value = "  +5.01%  "

# Row 6 - Solana
Set-TextValue "D6" "109.41"
$ws.Range("E6").Value = "  -0.65%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +2.88%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.12%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.642"
$ws.Range("E9").Value = "  +6.26%  "

# Row 10 - Avalanche
Set-TextValue "D10" "43.06"
$ws.Range("E10").Value = "  -2.06%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +3.34%  "

# Row 12 - Polkadot
Set-TextValue "D12" "8.85"
$ws.Range("E12").Value = "  +0.78%  "

# Row 13 - Polygon
$ws.Range("E13").Value = "  +3.89%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +2.15%  "

# Row 15 - Chainlink
Set-TextValue "D15" "16.39"
$ws.Range("E15").Value = "  +8.82%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue "D16" "2.715.00"
$ws.Range("E16").Value = "  +6.30%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.424.61"
$ws.Range("E17").Value = "  +8.62%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "43.385.00"
$ws.Range("E18").Value = "  +2.38%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  +3.64%  "

# Row 20 - Uniswap
Set-TextValue "D20" "7.25"
$ws.Range("E20").Value = "  -2.05%  "

# Row 21 - Litecoin
Set-TextValue "D21" "75.16"
$ws.Range("E21").Value = "  +3.99%  "

# Row 22 - PancakeSwap
Set-TextValue "D22" "3.43"
$ws.Range("E22").Value = "  -1.30%  "

# Row 23 - ImmutableX
$ws.Range("E23").Value = "  +9.34%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "258.45"
$ws.Range("E24").Value = "  +12.91%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("E25").Value = "  +0.45%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  +3.55%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.04%  "

# Row 28 - InjectiveProtocol
Set-TextValue "D28" "39.18"
$ws.Range("E28").Value = "  +2.77%  "

# Row 29 - Toncoin
Set-TextValue "D29" "2.25"
$ws.Range("E29").Value = "  +0.97%  "

# Row 30 - EthereumClassic
$ws.Range("E30").Value = "  +7.49%  "

# Row 31 - WEMIXToken
$ws.Range("E31").Value = "  -0.22%  "

# Row 32 - Monero
Set-TextValue "D32" "173.24"
$ws.Range("E32").Value = "  -0.30%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.0929"
$ws.Range("E33").Value = "  +3.45%  "

# Row 34 - Filecoin
Set-TextValue "D34" "6.02"
$ws.Range("E34").Value = "  +5.68%  "

# Row 35 - Stellar
$ws.Range("E35").Value = "  +5.52%  "

# Row 36 - RenderToken
Set-TextValue "D36" "4.97"
$ws.Range("E36").Value = "  -3.31%  "

# Row 37 - NEARProtocol
Set-TextValue "D37" "4.14"
$ws.Range("E37").Value = "  -4.98%  "

# Row 38 - VeChain
Set-TextValue "D38" "0.0376"
$ws.Range("E38").Value = "  -0.14%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -0.13%  "

# Row 40 - LidoDAOToken
Set-TextValue "D40" "2.81"
$ws.Range("E40").Value = "  +16.06%  "

# Row 41 - ARBITRUM
$ws.Range("E41").Value = "  +14.27%  "

# Row 42 - MultiversX
Set-TextValue "D42" "72.05"
$ws.Range("E42").Value = "  +0.20%  "

# Row 43 - Algorand
$ws.Range("E43").Value = "  -1.01%  "

# Rows 44 & 45 swapped: Celestia <-> FirstDigitalUSD
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D44" "1.00"
$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D45" "12.75"
$ws.Range("E45").Value = "  +1.02%  "

# Row 46 - THORChain
Set-TextValue "D46" "5.63"
$ws.Range("E46").Value = "  +3.62%  "

# Row 47 - FraxShare
Set-TextValue "D47" "9.38"
$ws.Range("E47").Value = "  +11.13%  "

# Row 48 - Aave
Set-TextValue "D48" "111.61"
$ws.Range("E48").Value = "  +7.74%  "

# Row 49 - TrustWalletToken
$ws.Range("E49").Value = "  +0.38%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  +2.88%  "

# Row 51 - WOONetwork
$ws.Range("E51").Value = "  +7.93%  "
